$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("user_names")

$ws.Range("B10").Value = "AswathDamodaran"
$ws.Range("B11").Value = "cstewartcfa"
$ws.Range("B12").Value = "BobPisani"

$ws.Range("B13").Select()

$ws.Columns.Item(2).ColumnWidth = 15.25

